$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "2023" column (E) next to the existing "2018" column (D) ---

# Header year cell E4 (copy formatting from D4, then set the year value)
$ws.Range("D4").Copy()
$ws.Range("E4").PasteSpecial(-4122)
$ws.Range("E4").Value = 2023

# Urban row value E5 (copy formatting from D5, then set the new figure)
$ws.Range("D5").Copy()
$ws.Range("E5").PasteSpecial(-4122)
$ws.Range("E5").Value = 38

# "Urbanisation" section header row 6 - add E6 matching D6's look, then re-bold both
$ws.Range("D6").Copy()
$ws.Range("E6").PasteSpecial(-4122)
$ws.Range("D6:E6").Font.Bold = $true

# Urban data row 7 - no 2023 figure available, show a dash, right aligned like D7
$ws.Range("D7").Copy()
$ws.Range("E7").PasteSpecial(-4122)
$ws.Range("E7").Value = "-"
$ws.Range("E7").HorizontalAlignment = -4152

# Rural data row 8 - no 2023 figure available, show a dash, right aligned like D8
$ws.Range("D8").Copy()
$ws.Range("E8").PasteSpecial(-4122)
$ws.Range("E8").Value = "-"
$ws.Range("E8").HorizontalAlignment = -4152

# --- Update the source footnote (row 9) to mention the new 2023 survey round ---
$ws.Range("A9").Value = " Көп көрсөткүчтүү кластердик изилдөөнүн маалыматтары боюнча, 2018-ж., 2023-ж."
$ws.Range("B9").Value = "По данным кластерного обследования по многим показателям, 2018г., 2023г."
$ws.Range("C9").Value = "According to the cluster survey in many respects, 2018, 2023."
